$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM_PURCHASE-EPS_V1")

$designators = "B1, BAT1, BAT2, BAT3, BAT4, C1, C2, C3, C4, C5, C6, C7, C8, C9, C10, C11, C12, C13, C14, C15, C16, C17, C18, C19, C20, C21, C22, C23, C24, C27, C28, C29, C30, C31, C32, C33, C34, C35, C36, C37, C38, C39, C40, C41, C42, C43, C44, C45, C46, C47, C48, C49, C50, C51, C52, C53, C54, C55, C56, C57, C58, C59, C60, C61, C62, C63, C64, C65, C68, C69, C70, C71, C72, C73, C74, C75, C76, C77, C78, C79, C80, C81, C82, C83, C84, C86, C87, C88, C89, C90, C91, C92, C93, C94, C95, FLT1, G1, L1, L2, L3, L4, LED1, R1, R2, R3, R4, R5, R6, R7, R8, R9, R10, R11, R12, R13, R14, R15, R16, R17, R18, R19, R20, R21, R22, R23, R24, R25, R26, R27, R28, R29, R30, R31, R32, R33, R34, R35, R36, R37, R38, R39, R40, R41, R42, R43, R44, R45, R46, R47, R48, R49, R50, R51, R52, R53, R54, R55, R56, R57, R58, R59, R60, R61, R62, R63, R64, R65, R66, R67, R68, R69, R70, R71, R72, R73, R74, R75, R76, R77, R78, R79, R80, R81, R82, R83, R84, R85, R86, R87, R88, R89, R90, R91, R92, R93, R94, R95, R96, R97, R98, R99, R100, R101, R102, R103, R104, R105, R106, R107, R108, R109, R110, R111, R112, R113, R114, R115, R116, R117, R118, R119, R121, R122, R123, R124, R125, R126, R127, R128, R129, R130, R131, R132, U1, U2, U7, U8, U9, U10, U11, U12, U13, U14, U15, U16, U17, U18, U19, U20, U21, U22, U23, U24, U25, USB1, VD1, VD2, VD3, VD4, VD5, VD9, VT1, VT2, VT3, VT4, VT5, VT6, VT7, VT8, VT9, VT10, VT11, VT12, VT13, VT14, VT15, VT16, VT17, VT18, VT19, X1, X2, X3, X4, X5, X6, X7, X8, X9, X10, X11, X12, X13, X14, X15"

$footprints = "button, BAT_18650_CLIPS, CERCAP 0402, CERCAP 0603, TECAP B, IND 0603, 3225, IND-5x5, IND-4.6x4.6, LED 0603, RES 0603, RES 0402, RES 0805, RES 1210, VSON-8, VSON-10, DFN-8-2x2, QFN-48, QFN-16(3,5x3,5), SO-8(EP), SOT223, SO-8, SOT-323-6 (SC-70), TO-92, USB_Type-C, SOD-123, SOT23-3, DFN-8_3x3, PicoBlade_4P, PC104_ArianSpace, PicoBlade_2P, PLS-2R-2.54, PBD-6-2.54"

$ws.Range("A3").Value = $designators
$ws.Range("C3").Value = 295
$ws.Range("D3").Value = $footprints
